# "picker updated with default value"
# Fill column G (the answer-explanation / picker column) for rows 26-100
# with the default placeholder text, matching the style already used by
# the other populated cells in that column (rows 1-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$placeholder = "I'm sure you know why (Placeholder)"

for ($row = 26; $row -le 100; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = $placeholder
    $cell.Font.Color = 3355443
}

# Update the view state: scroll position and active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 87
$win.ScrollColumn = 1
$ws.Range("G105").Select()
